# Updates cryptos list values (Coin/Link/Price/Volume) to match latest scrape.
# Price strings that look numeric (e.g. "624.60", "0.0000261") are written with
# a text NumberFormat first so Excel stores them verbatim instead of coercing them
# to floating point numbers (which would corrupt trailing zeros / precision).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Cells.Item(2, 4).Value2 = '82.162.11'
$ws.Cells.Item(2, 5).Value2 = '  +2.85%  '

# Row 3 (Ethereum)
$ws.Cells.Item(3, 4).Value2 = '3.197.53'
$ws.Cells.Item(3, 5).Value2 = '  -0.23%  '

# Row 4 (TetherUSD)
$ws.Cells.Item(4, 5).Value2 = '  +0.10%  '

# Row 5 (Solana)
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value2 = '216.59'
$ws.Cells.Item(5, 5).Value2 = '  +5.08%  '

# Row 6 (BNB)
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value2 = '624.60'
$ws.Cells.Item(6, 5).Value2 = '  -1.85%  '

# Row 7 (Dogecoin)
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value2 = '0.294'
$ws.Cells.Item(7, 5).Value2 = '  +22.64%  '

# Row 8 (USDC)
$ws.Cells.Item(8, 5).Value2 = '  +0.01%  '

# Row 9 (XRP)
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value2 = '0.587'
$ws.Cells.Item(9, 5).Value2 = '  +0.49%  '

# Row 10 (LidoStakedEther)
$ws.Cells.Item(10, 4).Value2 = '3.194.45'
$ws.Cells.Item(10, 5).Value2 = '  -0.22%  '

# Row 11 (Cardano)
$ws.Cells.Item(11, 5).Value2 = '  +2.00%  '

# Row 12 (ShibaInu)
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value2 = '0.0000261'
$ws.Cells.Item(12, 5).Value2 = '  +12.68%  '

# Row 14 (WrappedliquidstakedEther2.0)
$ws.Cells.Item(14, 2).Value2 = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(14, 3).Value2 = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(14, 4).Value2 = '3.791.34'
$ws.Cells.Item(14, 5).Value2 = '  +0.14%  '

# Row 15 (Toncoin)
$ws.Cells.Item(15, 2).Value2 = 'Toncoin'
$ws.Cells.Item(15, 3).Value2 = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value2 = '5.33'
$ws.Cells.Item(15, 5).Value2 = '  -3.73%  '

# Row 16 (Avalanche)
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value2 = '32.03'
$ws.Cells.Item(16, 5).Value2 = '  +0.30%  '

# Row 17 (WrappedBTC)
$ws.Cells.Item(17, 4).Value2 = '81.903.90'
$ws.Cells.Item(17, 5).Value2 = '  +2.97%  '

# Row 18 (WrappedEther)
$ws.Cells.Item(18, 4).Value2 = '3.197.67'
$ws.Cells.Item(18, 5).Value2 = '  +0.33%  '

# Row 19 (SuiNetwork)
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value2 = '3.28'
$ws.Cells.Item(19, 5).Value2 = '  +7.99%  '

# Row 20 (Chainlink)
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value2 = '14.15'
$ws.Cells.Item(20, 5).Value2 = '  -2.53%  '

# Row 21 (BitcoinCash)
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value2 = '437.33'
$ws.Cells.Item(21, 5).Value2 = '  +0.89%  '

# Row 22 (Uniswap)
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value2 = '9.03'
$ws.Cells.Item(22, 5).Value2 = '  -2.07%  '

# Row 23 (Polkadot)
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value2 = '5.14'
$ws.Cells.Item(23, 5).Value2 = '  -0.10%  '

# Row 24 (LEO)
$ws.Cells.Item(24, 5).Value2 = '  +6.04%  '

# Row 25 (NEARProtocol)
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value2 = '5.38'
$ws.Cells.Item(25, 5).Value2 = '  +13.32%  '

# Row 26 (WrappedeETH)
$ws.Cells.Item(26, 4).Value2 = '3.364.65'
$ws.Cells.Item(26, 5).Value2 = '  -0.03%  '

# Row 27 (Aptos)
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value2 = '11.21'
$ws.Cells.Item(27, 5).Value2 = '  -0.64%  '

# Row 28 (Litecoin)
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value2 = '76.80'
$ws.Cells.Item(28, 5).Value2 = '  -0.38%  '

# Row 29 (Dai)
$ws.Cells.Item(29, 5).Value2 = '  +0.44%  '

# Row 30 (PEPE)
$ws.Cells.Item(30, 5).Value2 = '  +3.97%  '

# Row 31 (Bittensor)
$ws.Cells.Item(31, 2).Value2 = 'Bittensor'
$ws.Cells.Item(31, 3).Value2 = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value2 = '588.33'
$ws.Cells.Item(31, 5).Value2 = '  +11.79%  '

# Row 32 (InternetComputer(DFINITY))
$ws.Cells.Item(32, 2).Value2 = 'InternetComputer(DFINITY)'
$ws.Cells.Item(32, 3).Value2 = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value2 = '9.14'
$ws.Cells.Item(32, 5).Value2 = '  +0.65%  '

# Row 33 (Binance-PegBSC-USD)
$ws.Cells.Item(33, 5).Value2 = '  -0.07%  '

# Row 34 (Fetch.AI)
$ws.Cells.Item(34, 5).Value2 = '  +1.55%  '

# Row 35 (Cronos)
$ws.Cells.Item(35, 2).Value2 = 'Cronos'
$ws.Cells.Item(35, 3).Value2 = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value2 = '0.144'
$ws.Cells.Item(35, 5).Value2 = '  +20.43%  '

# Row 36 (Kaspa)
$ws.Cells.Item(36, 2).Value2 = 'Kaspa'
$ws.Cells.Item(36, 3).Value2 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value2 = '0.155'
$ws.Cells.Item(36, 5).Value2 = '  +9.53%  '

# Row 37 (PancakeSwap)
$ws.Cells.Item(37, 5).Value2 = '  +0.80%  '

# Row 38 (EthereumClassic)
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value2 = '22.89'
$ws.Cells.Item(38, 5).Value2 = '  -1.26%  '

# Row 39 (FirstDigitalUSD)
$ws.Cells.Item(39, 2).Value2 = 'FirstDigitalUSD'
$ws.Cells.Item(39, 3).Value2 = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value2 = '1.00'
$ws.Cells.Item(39, 5).Value2 = '  +0.06%  '

# Row 40 (RenderToken)
$ws.Cells.Item(40, 2).Value2 = 'RenderToken'
$ws.Cells.Item(40, 3).Value2 = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value2 = '6.16'
$ws.Cells.Item(40, 5).Value2 = '  +10.97%  '

# Row 41 (PolygonEcosystemToken)
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value2 = '0.410'
$ws.Cells.Item(41, 5).Value2 = '  +0.23%  '

# Row 42 (Stacks)
$ws.Cells.Item(42, 5).Value2 = '  +13.90%  '

# Row 43 (dogwifhat)
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value2 = '3.08'
$ws.Cells.Item(43, 5).Value2 = '  +21.24%  '

# Row 44 (WhiteBITCoin)
$ws.Cells.Item(44, 5).Value2 = '  +3.89%  '

# Row 45 (Monero)
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value2 = '161.26'
$ws.Cells.Item(45, 5).Value2 = '  -2.48%  '

# Row 46 (USDe)
$ws.Cells.Item(46, 5).Value2 = '  +0.03%  '

# Row 47 (Aave)
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value2 = '188.89'
$ws.Cells.Item(47, 5).Value2 = '  -2.04%  '

# Row 48 (OKB)
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value2 = '44.78'
$ws.Cells.Item(48, 5).Value2 = '  +3.54%  '

# Row 49 (ImmutableX)
$ws.Cells.Item(49, 5).Value2 = '  +1.17%  '

# Row 50 (Mantle)
$ws.Cells.Item(50, 2).Value2 = 'Mantle'
$ws.Cells.Item(50, 3).Value2 = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value2 = '0.776'
$ws.Cells.Item(50, 5).Value2 = '  -5.58%  '

# Row 51 (InjectiveProtocol)
$ws.Cells.Item(51, 2).Value2 = 'InjectiveProtocol'
$ws.Cells.Item(51, 3).Value2 = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value2 = '26.34'
$ws.Cells.Item(51, 5).Value2 = '  +1.62%  '
